$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.870.36"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "2.222.39"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -1.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.95"
$ws.Range("E5").Value = "  -3.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.23"
$ws.Range("E6").Value = "  -4.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  -3.85%  "
$ws.Range("E9").Value = "  -6.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.91"
$ws.Range("E10").Value = "  -5.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0775"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.93"
$ws.Range("E12").Value = "  -4.18%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "2.560.30"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "2.217.00"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.44"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("E17").Value = "  -7.47%  "
$ws.Range("D18").Value = "43.753.89"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "0.0₃0899"
$ws.Range("E19").Value = "  -6.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.92"
$ws.Range("E21").Value = "  -6.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.46"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.37"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("E24").Value = "  -5.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("E26").Value = "  -6.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.22"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  -5.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "152.38"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.18"
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.39"
$ws.Range("E32").Value = "  -9.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0754"
$ws.Range("E33").Value = "  -5.92%  "
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.84"
$ws.Range("E36").Value = "  -8.77%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.102"
$ws.Range("E37").Value = "  -6.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.68"
$ws.Range("E38").Value = "  -5.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0299"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -4.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.17"
$ws.Range("E41").Value = "  -7.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.13"
$ws.Range("E42").Value = "  -8.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "1.828.41"
$ws.Range("E44").Value = "  +3.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.77"
$ws.Range("E45").Value = "  +11.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.181"
$ws.Range("E46").Value = "  -6.15%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.06"
$ws.Range("E47").Value = "  -5.23%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "66.90"
$ws.Range("E48").Value = "  -4.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.88"
$ws.Range("E49").Value = "  -9.17%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.442.11"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "13.72"
$ws.Range("E51").Value = "  -4.55%  "

Write-Output "done"
